$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing D:K data to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats from the (now shifted) F:G columns into the new D:E columns
# so the new columns inherit the correct date/number styles row by row.
$ws.Range("F7:G102").Copy() | Out-Null
$ws.Range("D7:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 43465
$arr[0,1] = 43373
$arr[0,2] = 43281
$arr[0,3] = 43190
$arr[0,4] = 43100
$arr[0,5] = 43008
$arr[0,6] = 42916
$arr[0,7] = 42825
$arr[0,8] = 42735
$arr[0,9] = 42643
$ws.Range("D7:M7").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 535700
$arr[0,1] = 538300
$arr[0,2] = 509100
$arr[0,3] = 366800
$arr[0,4] = 410300
$arr[0,5] = 404000
$arr[0,6] = 384900
$arr[0,7] = 289200
$arr[0,8] = 1104100
$arr[0,9] = 297100
$ws.Range("D8:M8").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 67900
$arr[0,1] = 56500
$arr[0,2] = 50700
$arr[0,3] = 38300
$arr[0,4] = 39000
$arr[0,5] = 33800
$arr[0,6] = 34900
$arr[0,7] = 29000
$arr[0,8] = 102900
$arr[0,9] = 25700
$ws.Range("D9:M9").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 467700
$arr[0,1] = 481800
$arr[0,2] = 458400
$arr[0,3] = 328500
$arr[0,4] = 371300
$arr[0,5] = 370200
$arr[0,6] = 350000
$arr[0,7] = 260100
$arr[0,8] = 1001300
$arr[0,9] = 271400
$ws.Range("D10:M10").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 73200
$arr[0,1] = 69600
$arr[0,2] = 58500
$arr[0,3] = 51400
$arr[0,4] = 53100
$arr[0,5] = 55300
$arr[0,6] = 48000
$arr[0,7] = 45700
$arr[0,8] = 161100
$arr[0,9] = 44400
$ws.Range("D12:M12").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D13:M13").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = "NA"
$arr[0,1] = "NA"
$arr[0,2] = "NA"
$arr[0,3] = "NA"
$arr[0,4] = 0
$arr[0,5] = "NA"
$arr[0,6] = "NA"
$arr[0,7] = "NA"
$arr[0,8] = -11600
$arr[0,9] = 0
$ws.Range("D14:M14").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D15:M15").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 430800
$arr[0,1] = 438700
$arr[0,2] = 398900
$arr[0,3] = 327200
$arr[0,4] = 319200
$arr[0,5] = 326500
$arr[0,6] = 298400
$arr[0,7] = 277800
$arr[0,8] = 1058600
$arr[0,9] = 276200
$ws.Range("D17:M17").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 104800
$arr[0,1] = 99500
$arr[0,2] = 110300
$arr[0,3] = 39600
$arr[0,4] = 91100
$arr[0,5] = 77500
$arr[0,6] = 86400
$arr[0,7] = 11300
$arr[0,8] = 45600
$arr[0,9] = 20800
$ws.Range("D18:M18").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -24400
$arr[0,1] = 15400
$arr[0,2] = 14000
$arr[0,3] = 1300
$arr[0,4] = -18400
$arr[0,5] = -12100
$arr[0,6] = 6200
$arr[0,7] = -14100
$arr[0,8] = -165400
$arr[0,9] = -50900
$ws.Range("D20:M20").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = "NA"
$arr[0,1] = "NA"
$arr[0,2] = "NA"
$arr[0,3] = "NA"
$arr[0,4] = "NA"
$arr[0,5] = "NA"
$arr[0,6] = "NA"
$arr[0,7] = "NA"
$arr[0,8] = "NA"
$arr[0,9] = "NA"
$ws.Range("D21:M21").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D22:M22").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 80400
$arr[0,1] = 114900
$arr[0,2] = 124200
$arr[0,3] = 40900
$arr[0,4] = 72700
$arr[0,5] = 65400
$arr[0,6] = 92600
$arr[0,7] = -2700
$arr[0,8] = -119800
$arr[0,9] = -30100
$ws.Range("D23:M23").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 13900
$arr[0,1] = 2500
$arr[0,2] = 17700
$arr[0,3] = 10400
$arr[0,4] = 6300
$arr[0,5] = 7600
$arr[0,6] = 8200
$arr[0,7] = -300
$arr[0,8] = -7400
$arr[0,9] = -1800
$ws.Range("D24:M24").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D25:M25").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 66500
$arr[0,1] = 112400
$arr[0,2] = 106500
$arr[0,3] = 30500
$arr[0,4] = 66400
$arr[0,5] = 57900
$arr[0,6] = 84400
$arr[0,7] = -2500
$arr[0,8] = -112400
$arr[0,9] = -28300
$ws.Range("D26:M26").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 61000
$arr[0,1] = 107800
$arr[0,2] = 101700
$arr[0,3] = 25900
$arr[0,4] = 61500
$arr[0,5] = 52400
$arr[0,6] = 80000
$arr[0,7] = -3100
$arr[0,8] = -114000
$arr[0,9] = -29000
$ws.Range("D27:M27").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D28:M28").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D29:M29").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D30:M30").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D31:M31").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 24400
$arr[0,1] = -15400
$arr[0,2] = -14000
$arr[0,3] = -1300
$arr[0,4] = 18400
$arr[0,5] = 12100
$arr[0,6] = -6200
$arr[0,7] = 14100
$arr[0,8] = 165400
$arr[0,9] = 50900
$ws.Range("D32:M32").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 61000
$arr[0,1] = 107800
$arr[0,2] = 101700
$arr[0,3] = 25900
$arr[0,4] = 61500
$arr[0,5] = 52400
$arr[0,6] = 80000
$arr[0,7] = -3100
$arr[0,8] = -114000
$arr[0,9] = -29000
$ws.Range("D33:M33").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D34:M34").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 61000
$arr[0,1] = 107800
$arr[0,2] = 101700
$arr[0,3] = 25900
$arr[0,4] = 61500
$arr[0,5] = 52400
$arr[0,6] = 80000
$arr[0,7] = -3100
$arr[0,8] = -114000
$arr[0,9] = -29000
$ws.Range("D35:M35").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 43465
$arr[0,1] = 43373
$arr[0,2] = 43281
$arr[0,3] = 43190
$arr[0,4] = 43100
$arr[0,5] = 43008
$arr[0,6] = 42916
$arr[0,7] = 42825
$arr[0,8] = 42735
$arr[0,9] = 42643
$ws.Range("D38:M38").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 354300
$arr[0,1] = 212100
$arr[0,2] = 280100
$arr[0,3] = 221600
$arr[0,4] = 226300
$arr[0,5] = 272200
$arr[0,6] = 385000
$arr[0,7] = 173500
$arr[0,8] = 174600
$arr[0,9] = 22100
$ws.Range("D41:M41").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 680800
$arr[0,1] = 730400
$arr[0,2] = 710400
$arr[0,3] = 611400
$arr[0,4] = 510200
$arr[0,5] = 372900
$arr[0,6] = 168300
$arr[0,7] = 169100
$arr[0,8] = 125000
$arr[0,9] = 21700
$ws.Range("D42:M42").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 136200
$arr[0,1] = 134500
$arr[0,2] = 103900
$arr[0,3] = 97900
$arr[0,4] = 99100
$arr[0,5] = 84500
$arr[0,6] = 78800
$arr[0,7] = 71800
$arr[0,8] = 61800
$arr[0,9] = 11300
$ws.Range("D43:M43").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D44:M44").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 241200
$arr[0,1] = 266200
$arr[0,2] = 270900
$arr[0,3] = 131900
$arr[0,4] = 111400
$arr[0,5] = 210100
$arr[0,6] = 269700
$arr[0,7] = 237200
$arr[0,8] = 229500
$arr[0,9] = 23000
$ws.Range("D45:M45").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1412600
$arr[0,1] = 1343200
$arr[0,2] = 1365200
$arr[0,3] = 1062800
$arr[0,4] = 947000
$arr[0,5] = 939800
$arr[0,6] = 901800
$arr[0,7] = 651600
$arr[0,8] = 590900
$arr[0,9] = 78100
$ws.Range("D46:M46").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 499500
$arr[0,1] = 489500
$arr[0,2] = 283300
$arr[0,3] = 268700
$arr[0,4] = 268400
$arr[0,5] = 297400
$arr[0,6] = 318000
$arr[0,7] = 307900
$arr[0,8] = 308100
$arr[0,9] = 54000
$ws.Range("D47:M47").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 197300
$arr[0,1] = 196500
$arr[0,2] = 194200
$arr[0,3] = 196200
$arr[0,4] = 200600
$arr[0,5] = 204300
$arr[0,6] = 209700
$arr[0,7] = 208700
$arr[0,8] = 215400
$arr[0,9] = 28800
$ws.Range("D48:M48").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 2519700
$arr[0,1] = 2526500
$arr[0,2] = 2532800
$arr[0,3] = 2541100
$arr[0,4] = 2549400
$arr[0,5] = 2563400
$arr[0,6] = 2571600
$arr[0,7] = 2528000
$arr[0,8] = 2536300
$arr[0,9] = 381100
$ws.Range("D49:M49").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D50:M50").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D51:M51").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 94900
$arr[0,1] = 81000
$arr[0,2] = 102800
$arr[0,3] = 223400
$arr[0,4] = 229600
$arr[0,5] = 163400
$arr[0,6] = 169600
$arr[0,7] = 17500
$arr[0,8] = 32500
$arr[0,9] = 7200
$ws.Range("D52:M52").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D53:M53").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 4724000
$arr[0,1] = 4636600
$arr[0,2] = 4478400
$arr[0,3] = 4292300
$arr[0,4] = 4195000
$arr[0,5] = 4168200
$arr[0,6] = 4170700
$arr[0,7] = 3713700
$arr[0,8] = 3683200
$arr[0,9] = 549200
$ws.Range("D54:M54").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 131700
$arr[0,1] = 147000
$arr[0,2] = 140200
$arr[0,3] = 123500
$arr[0,4] = 92700
$arr[0,5] = 89800
$arr[0,6] = 100300
$arr[0,7] = 88500
$arr[0,8] = 89000
$arr[0,9] = 11900
$ws.Range("D57:M57").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 120600
$arr[0,1] = 120900
$arr[0,2] = 116700
$arr[0,3] = 11100
$arr[0,4] = 11100
$arr[0,5] = 86300
$arr[0,6] = 158500
$arr[0,7] = 266000
$arr[0,8] = 268000
$arr[0,9] = 38600
$ws.Range("D58:M58").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 828800
$arr[0,1] = 816000
$arr[0,2] = 809800
$arr[0,3] = 753600
$arr[0,4] = 727500
$arr[0,5] = 710800
$arr[0,6] = 688300
$arr[0,7] = 670700
$arr[0,8] = 643600
$arr[0,9] = 96100
$ws.Range("D59:M59").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1081200
$arr[0,1] = 1083900
$arr[0,2] = 1066700
$arr[0,3] = 888200
$arr[0,4] = 831200
$arr[0,5] = 886900
$arr[0,6] = 947100
$arr[0,7] = 1025200
$arr[0,8] = 1000600
$arr[0,9] = 146700
$ws.Range("D60:M60").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 5600
$arr[0,2] = 5600
$arr[0,3] = 111500
$arr[0,4] = 115400
$arr[0,5] = 128100
$arr[0,6] = 130300
$arr[0,7] = 21800
$arr[0,8] = 21800
$arr[0,9] = 0
$ws.Range("D61:M61").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 42300
$arr[0,1] = 43900
$arr[0,2] = 44000
$arr[0,3] = 47000
$arr[0,4] = 50000
$arr[0,5] = 52900
$arr[0,6] = 55900
$arr[0,7] = 57700
$arr[0,8] = 64500
$arr[0,9] = 9800
$ws.Range("D62:M62").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D63:M63").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D64:M64").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D65:M65").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1455900
$arr[0,1] = 1449300
$arr[0,2] = 1415400
$arr[0,3] = 1325400
$arr[0,4] = 1287200
$arr[0,5] = 1335200
$arr[0,6] = 1412300
$arr[0,7] = 1127700
$arr[0,8] = 1109300
$arr[0,9] = 159600
$ws.Range("D66:M66").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D68:M68").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D69:M69").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D70:M70").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D71:M71").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 65200
$arr[0,1] = -900
$arr[0,2] = -113700
$arr[0,3] = -220100
$arr[0,4] = -250800
$arr[0,5] = -317100
$arr[0,6] = -374200
$arr[0,7] = -449100
$arr[0,8] = -446600
$arr[0,9] = -67200
$ws.Range("D72:M72").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D73:M73").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D74:M74").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D75:M75").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 3268100
$arr[0,1] = 3187300
$arr[0,2] = 3063000
$arr[0,3] = 2966800
$arr[0,4] = 2907900
$arr[0,5] = 2833000
$arr[0,6] = 2758500
$arr[0,7] = 2586000
$arr[0,8] = 2573900
$arr[0,9] = 389500
$ws.Range("D76:M76").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D77:M77").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 43465
$arr[0,1] = 43373
$arr[0,2] = 43281
$arr[0,3] = 43190
$arr[0,4] = 43100
$arr[0,5] = 43008
$arr[0,6] = 42916
$arr[0,7] = 42825
$arr[0,8] = 42735
$arr[0,9] = 42643
$ws.Range("D80:M80").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 61000
$arr[0,1] = 107800
$arr[0,2] = 101700
$arr[0,3] = 25900
$arr[0,4] = 61500
$arr[0,5] = 52400
$arr[0,6] = 80000
$arr[0,7] = -3100
$arr[0,8] = -114000
$arr[0,9] = -29000
$ws.Range("D81:M81").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D83:M83").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D84:M84").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D85:M85").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D86:M86").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D87:M87").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D88:M88").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D89:M89").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D91:M91").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D92:M92").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D93:M93").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D94:M94").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D96:M96").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D97:M97").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D98:M98").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D99:M99").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D100:M100").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D101:M101").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D102:M102").Value = $arr
